$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns (D/E) store values as text in the source data
# (e.g. "67.460.06", "0.999", "  +0.81%  "). Plain Value assignment would let
# Excel auto-coerce purely numeric-looking strings (e.g. "0.999") into real
# numbers, so for those specific cells we force a text number format first.

$ws.Range("D2").Value = '67.460.06'
$ws.Range("E2").Value = '  +0.81%  '

$ws.Range("D3").Value = '3.494.87'
$ws.Range("E3").Value = '  -0.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.03'
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.98'
$ws.Range("E6").Value = '  +4.29%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  +2.09%  '

$ws.Range("D9").Value = '3.496.33'
$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("E10").Value = '  +5.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.05'
$ws.Range("E11").Value = '  -2.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.436'
$ws.Range("E12").Value = '  +0.98%  '

$ws.Range("D13").Value = '4.087.30'
$ws.Range("E13").Value = '  -0.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.44'
$ws.Range("E14").Value = '  +11.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.136'
$ws.Range("E15").Value = '  +1.03%  '

$ws.Range("D16").Value = '67.410.57'
$ws.Range("E16").Value = '  +0.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000178'
$ws.Range("E17").Value = '  -0.48%  '

$ws.Range("D18").Value = '3.479.73'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.29'
$ws.Range("E19").Value = '  +0.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.32'
$ws.Range("E20").Value = '  +0.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.25'
$ws.Range("E21").Value = '  -1.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.97'
$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.07'
$ws.Range("E23").Value = '  +0.81%  '

$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.542'
$ws.Range("E24").Value = '  +1.07%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.75'
$ws.Range("E26").Value = '  +0.83%  '

$ws.Range("E27").Value = '  +0.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.39'
$ws.Range("E28").Value = '  +1.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.175'
$ws.Range("E29").Value = '  -3.52%  '

$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.23'
$ws.Range("E31").Value = '  +0.79%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.43'
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.07'
$ws.Range("E33").Value = '  +0.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.56'
$ws.Range("E34").Value = '  -0.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.40'
$ws.Range("E35").Value = '  +0.48%  '

$ws.Range("E36").Value = '  +0.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.60'
$ws.Range("E37").Value = '  -0.73%  '

$ws.Range("E38").Value = '  +1.12%  '

$ws.Range("E39").Value = '  -1.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.80'
$ws.Range("E40").Value = '  +10.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.89'
$ws.Range("E41").Value = '  -1.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.80'
$ws.Range("E42").Value = '  -1.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.66'
$ws.Range("E43").Value = '  +0.44%  '

$ws.Range("D44").Value = '2.851.48'
$ws.Range("E44").Value = '  +0.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '27.10'
$ws.Range("E45").Value = '  +0.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.26'
$ws.Range("E46").Value = '  +0.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0725'
$ws.Range("E47").Value = '  -1.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '41.75'
$ws.Range("E48").Value = '  -2.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0302'
$ws.Range("E49").Value = '  -0.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '336.17'
$ws.Range("E50").Value = '  -0.25%  '

$ws.Range("E51").Value = '  -1.96%  '
